# Update the division-problem values in the practice-sheet table.
# The worksheet table has 5 "data" rows (each followed by 3 blank rows)
# and 5 columns; only the data-row cells carry text. Target cells are
# addressed directly by (row, column) via Table.Cell(row, col) so that
# duplicate old values (e.g. "66÷8=" appears twice but maps to two
# different replacements) are updated independently and correctly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=1;  Col=1; Old="95÷8="; New="47÷3="},
    @{Row=1;  Col=2; Old="62÷2="; New="53÷3="},
    @{Row=1;  Col=3; Old="52÷8="; New="29÷5="},
    @{Row=1;  Col=4; Old="69÷5="; New="87÷6="},
    @{Row=1;  Col=5; Old="63÷6="; New="23÷9="},

    @{Row=5;  Col=1; Old="91÷2="; New="95÷7="},
    @{Row=5;  Col=2; Old="99÷9="; New="75÷2="},
    @{Row=5;  Col=3; Old="79÷7="; New="88÷8="},
    @{Row=5;  Col=4; Old="55÷2="; New="79÷6="},
    @{Row=5;  Col=5; Old="87÷4="; New="32÷7="},

    @{Row=9;  Col=1; Old="74÷4="; New="83÷9="},
    @{Row=9;  Col=2; Old="66÷8="; New="33÷5="},
    @{Row=9;  Col=3; Old="79÷2="; New="91÷5="},
    @{Row=9;  Col=4; Old="12÷5="; New="23÷4="},
    @{Row=9;  Col=5; Old="28÷6="; New="38÷8="},

    @{Row=13; Col=1; Old="66÷8="; New="90÷6="},
    @{Row=13; Col=2; Old="35÷7="; New="60÷7="},
    @{Row=13; Col=3; Old="84÷7="; New="91÷3="},
    @{Row=13; Col=4; Old="56÷4="; New="37÷6="},
    @{Row=13; Col=5; Old="75÷7="; New="94÷4="},

    @{Row=17; Col=1; Old="12÷8="; New="39÷4="},
    @{Row=17; Col=2; Old="35÷2="; New="59÷7="},
    @{Row=17; Col=3; Old="26÷7="; New="81÷8="},
    @{Row=17; Col=4; Old="44÷7="; New="64÷9="},
    @{Row=17; Col=5; Old="20÷2="; New="35÷7="}
)

foreach ($chg in $changes) {
    $cell = $t.Cell($chg.Row, $chg.Col)
    $rng = $cell.Range
    # Sanity-check we are about to overwrite the expected old value before
    # mutating (guards against any cell mis-addressing).
    if ($rng.Text -notlike ($chg.Old + "*")) {
        Write-Host "WARNING: cell ($($chg.Row),$($chg.Col)) expected '$($chg.Old)' but found '$($rng.Text)'"
    }
    $rng.Text = $chg.New
}

Write-Host "Applied $($changes.Count) replacements"
